$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Feature branch"
$ws.Range("A2").Select() | Out-Null
